$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

# ---------------------------------------------------------------------------
# 1. Status text: "Ready for handoff" -> "Handed back: in sync with en-US"
#    This literal string is used independently by each cell, so every cell
#    that shows it needs to be rewritten explicitly.
# ---------------------------------------------------------------------------
$newStatus = "Handed back: in sync with en-US"

$wsOverview.Range("E2").Value = $newStatus
$wsOverview.Range("F2").Value = $newStatus
$wsOverview.Range("E3").Value = $newStatus
$wsOverview.Range("F3").Value = $newStatus

$wsZhCn.Range("C2").Value = $newStatus
$wsZhCn.Range("C3").Value = $newStatus

$wsDeDe.Range("C2").Value = $newStatus
$wsDeDe.Range("C3").Value = $newStatus

# ---------------------------------------------------------------------------
# 2. Populate the "Latest Target File" / "Latest Handback File" columns
#    (I and J) for both language sheets, turning column I into a hyperlink
#    pointing at the same markdown file as column A.
# ---------------------------------------------------------------------------

# zh-cn sheet
$wsZhCn.Range("I2").Value = "a7ef105e-fe44-4ab2-a245-936e4a263d00.md"
$wsZhCn.Range("I2").Style = "Hyperlink"
$wsZhCn.Range("J2").Value = "a7ef105e-fe44-4ab2-a245-936e4a263d00.b74257a9eabb71db64ec0b14d92c126457c901cd.zh-cn.xlf"

$wsZhCn.Range("I3").Value = "f68e5e0c-b2a4-47e5-8b0f-52d16917e7ba.md"
$wsZhCn.Range("I3").Style = "Hyperlink"
$wsZhCn.Range("J3").Value = "f68e5e0c-b2a4-47e5-8b0f-52d16917e7ba.8dc492555c187a36dba8c74eb8f5244b4f7934ae.zh-cn.xlf"

# de-de sheet
$wsDeDe.Range("I2").Value = "a7ef105e-fe44-4ab2-a245-936e4a263d00.md"
$wsDeDe.Range("I2").Style = "Hyperlink"
$wsDeDe.Range("J2").Value = "a7ef105e-fe44-4ab2-a245-936e4a263d00.b74257a9eabb71db64ec0b14d92c126457c901cd.de-de.xlf"
$wsDeDe.Range("K2").Value = "2016-08-17 22:46:36"

$wsDeDe.Range("I3").Value = "f68e5e0c-b2a4-47e5-8b0f-52d16917e7ba.md"
$wsDeDe.Range("I3").Style = "Hyperlink"
$wsDeDe.Range("J3").Value = "f68e5e0c-b2a4-47e5-8b0f-52d16917e7ba.8dc492555c187a36dba8c74eb8f5244b4f7934ae.de-de.xlf"
$wsDeDe.Range("K3").Value = "2016-08-17 22:46:36"

# ---------------------------------------------------------------------------
# 3. Rebuild the hyperlinks on zh-cn / de-de so the ordering/ids match:
#    A2, I2, A3, I3 (delete then re-add in that order).
# ---------------------------------------------------------------------------
$githubBase = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/e12470103b271a49d319fc31c15edc725f5f3f74/e2e/"

$wsZhCn.Hyperlinks.Delete()
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A2"), $githubBase + "a7ef105e-fe44-4ab2-a245-936e4a263d00.md", $null, $null, "a7ef105e-fe44-4ab2-a245-936e4a263d00.md")
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("I2"), $githubBase + "a7ef105e-fe44-4ab2-a245-936e4a263d00.md", $null, $null, "a7ef105e-fe44-4ab2-a245-936e4a263d00.md")
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A3"), $githubBase + "f68e5e0c-b2a4-47e5-8b0f-52d16917e7ba.md", $null, $null, "f68e5e0c-b2a4-47e5-8b0f-52d16917e7ba.md")
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("I3"), $githubBase + "f68e5e0c-b2a4-47e5-8b0f-52d16917e7ba.md", $null, $null, "f68e5e0c-b2a4-47e5-8b0f-52d16917e7ba.md")

$wsDeDe.Hyperlinks.Delete()
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A2"), $githubBase + "a7ef105e-fe44-4ab2-a245-936e4a263d00.md", $null, $null, "a7ef105e-fe44-4ab2-a245-936e4a263d00.md")
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("I2"), $githubBase + "a7ef105e-fe44-4ab2-a245-936e4a263d00.md", $null, $null, "a7ef105e-fe44-4ab2-a245-936e4a263d00.md")
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A3"), $githubBase + "f68e5e0c-b2a4-47e5-8b0f-52d16917e7ba.md", $null, $null, "f68e5e0c-b2a4-47e5-8b0f-52d16917e7ba.md")
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("I3"), $githubBase + "f68e5e0c-b2a4-47e5-8b0f-52d16917e7ba.md", $null, $null, "f68e5e0c-b2a4-47e5-8b0f-52d16917e7ba.md")

# ---------------------------------------------------------------------------
# 4. Column width adjustments.
#    ColumnWidth is expressed in "characters"; the engine stores
#    width = round(ColumnWidth*6)/6 + 5/6 in the workbook XML, so the inputs
#    below are chosen to land as closely as possible on the target widths.
# ---------------------------------------------------------------------------

# Overview: columns E (zh-cn) and F (de-de) grow from ~17.22 to ~29.98
$wsOverview.Columns.Item(5).ColumnWidth = 29.166666666666668
$wsOverview.Columns.Item(6).ColumnWidth = 29.166666666666668

# zh-cn / de-de: column C (Status) grows from ~17.22 to ~29.98
$wsZhCn.Columns.Item(3).ColumnWidth = 29.166666666666668
$wsDeDe.Columns.Item(3).ColumnWidth = 29.166666666666668

# zh-cn / de-de: columns I (Latest Target File) and J (Latest Handback File) -> 40
$wsZhCn.Columns.Item(9).ColumnWidth = 39.166666666666664
$wsZhCn.Columns.Item(10).ColumnWidth = 39.166666666666664
$wsDeDe.Columns.Item(9).ColumnWidth = 39.166666666666664
$wsDeDe.Columns.Item(10).ColumnWidth = 39.166666666666664
